$d = $word.ActiveDocument
$find = $d.Content.Find

# Title
$find.Execute(
    'Unraveling Gender Disparity in STEM',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'The Fascinating Realm of Chemistry: Unveiling the Secrets of Matter', 2
) | Out-Null

# Author name
$find.Execute(
    'Sarah Miller',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'Dennis Wilson', 2
) | Out-Null

# Author email
$find.Execute(
    'sarahmiller@gmail.com',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'dennis.wilson@newwave.edu', 2
) | Out-Null

# Body paragraph - block 1
$find.Execute(
    'Throughout history, the realm of science, technology, engineering, and mathematics (STEM) has been predominantly male-dominated. This disparity is a global phenomenon, transcending cultural and socioeconomic boundaries. Consequently, society has missed out on the invaluable contributions of women in STEM fields, resulting in a skewed representation of perspectives, ideas, and innovations. This essay delves into the multifaceted issue of gender disparity in STEM, examining its root causes, detrimental effects, and potential solutions to foster a more inclusive environment.',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'In the vast tapestry of scientific disciplines, Chemistry stands as a beacon of discovery, shedding light on the intricate workings of matter. It is a subject that touches every aspect of our lives, from the food we consume to the air we breathe, and even the clothes we wear. In this exploration, we will delve into the captivating world of Chemistry, revealing the fundamental concepts that govern the behavior of substances and the remarkable transformations they undergo.', 2
) | Out-Null

# Body paragraph - block 2
$find.Execute(
    'Gender stereotypes and societal expectations play a pivotal role in shaping girls'' and women''s career choices. From an early age, girls are often discouraged from pursuing careers in STEM, as these fields are traditionally perceived as masculine. This perception is reinforced by the limited visibility of female role models in STEM, perpetuating the notion that these fields are not suitable for women. Furthermore, the lack of gender diversity in STEM creates a hostile environment for women, where they face discrimination, prejudice, and a lack of support.',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'Within the realm of Chemistry, we will unravel the enigmatic nature of atoms and molecules, the building blocks of all matter. We will explore the periodic table, a roadmap that unveils the properties and relationships of these fundamental particles. Through engaging experiments and captivating demonstrations, we will witness the symphony of chemical reactions, marveling at the vibrant colors and intriguing changes that accompany them.', 2
) | Out-Null

# Body paragraph - block 3
$find.Execute(
    'The consequences of gender disparity in STEM are multi-faceted and far-reaching. It deprives society of the talents and contributions of a large pool of potential scientists, engineers, and innovators. This has implications for economic growth, as well as the development of new technologies and solutions to address global challenges. Moreover, the underrepresentation of women in STEM reinforces gender stereotypes and perpetuates the cycle of discrimination, creating a vicious loop that is difficult to break.',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'Furthermore, we will probe the depths of chemical bonding, the forces that hold atoms together and determine the properties of substances. We will delve into the mysteries of acids and bases, unveiling their roles in everyday phenomena and their applications in various industries. As we progress in our journey through Chemistry, we will appreciate the intricate dance of particles and the elegance of chemical principles that orchestrate the material world around us.', 2
) | Out-Null

# Summary paragraph
$find.Execute(
    'The gender disparity in STEM is a multifaceted issue with profound implications for society. Rooted in societal stereotypes and a lack of female role models, it results in a hostile environment for women in STEM, leading to discrimination, prejudice, and a lack of support. This disparity has detrimental consequences, including the loss of talent, the perpetuation of gender stereotypes, and the stifling of innovation. To address this issue, comprehensive efforts are required to challenge stereotypes, promote female role models, create inclusive environments, and implement policies that support women in STEM. By fostering a more diverse and inclusive STEM workforce, society can unlock the full potential of innovation and progress.',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'Chemistry, a captivating realm of scientific inquiry, uncovers the secrets of matter and its transformations. Through the study of atoms, molecules, and chemical reactions, we gain profound insights into the behavior of substances and their applications in various industries. By delving into the fundamental principles of Chemistry, we not only enhance our understanding of the world around us but also equip ourselves with valuable tools for solving real-world problems, paving the way for future scientific advancements and technological breakthroughs.', 2
) | Out-Null

# Append a new trailing empty paragraph at the very end of the document (before sectPr)
$endPos = $d.Content.End
$endRange = $d.Range($endPos, $endPos)
$endRange.Text = "`r"
